{"js": "// The document uses a Jinja-style merge field \"{{ doc_identidade_assistido }}\"\n// in two places (the intro declaration paragraph and the identity-document\n// line near the signature block). Both need to become\n// \"{{ sigla_identidade }}/{{ sigla_estado_identidade }}\".\n//\n// Rather than juggling the many tiny runs that make up the template markup,\n// search the body for the literal field-name text \"doc_identidade_assistido\"\n// (it only ever appears as that whole token, never as a substring of another\n// field name) and replace each hit in place. This preserves the surrounding\n// \"{{ \" / \" }}\" delimiters and any other text untouched.\n\nconst results = context.document.body.search(\"doc_identidade_assistido\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (const r of results.items) {\n  r.insertText(\"sigla_identidade }}/{{ sigla_estado_identidade\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The template merge field \"{{ doc_identidade_assistido }}\" is used twice in\n# this document (the intro declaration paragraph and the identity-document\n# line near the signature block). Both occurrences need to become\n# \"{{ sigla_identidade }}/{{ sigla_estado_identidade }}\".\n#\n# The field name \"doc_identidade_assistido\" only ever appears as that whole\n# token (never as a substring of another field name such as\n# \"num_identidade_assistido\" or \"nome_assistido\"), so a simple\n# Find/Replace over the whole document body for that literal text is safe\n# and precise, leaving the surrounding \"{{ \" / \" }}\" delimiters untouched.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"doc_identidade_assistido\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"sigla_identidade }}/{{ sigla_estado_identidade\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n"}
